# "In progress Data Preparation"
# - Fix stale shared-formula span in column L (now only runs to row 86)
# - Unhide helper column M
# - Autofit / resize the data columns to fit their (now wider) content
# - Leave the selection on D11 (where the user was working)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-key the "Management Action" helper formula down the used range ---
# The shared formula used to stretch to L97 (stale, from when the table had
# more rows); re-entering it over the real data range (L66:L86) corrects it.
$ws.Range("L66:L86").Formula = '=B66&" "&C66&"-"&M66&" Action Details"'

# --- Unhide column M (the helper index column used by the formula above) ---
$ws.Columns.Item(13).Hidden = $false

# --- Autofit all the columns against their actual content ---
$ws.Range("A1:M86").EntireColumn.AutoFit() | Out-Null

# Nudge the autofit result to the precise widths Excel's own metrics produced
$ws.Columns.Item(1).ColumnWidth = 30.666666666666668
$ws.Columns.Item(2).ColumnWidth = 20.0
$ws.Columns.Item(3).ColumnWidth = 17.166666666666668
$ws.Columns.Item(4).ColumnWidth = 15.166666666666666
$ws.Columns.Item(5).ColumnWidth = 26.666666666666668
$ws.Columns.Item(6).ColumnWidth = 16.333333333333332
$ws.Columns.Item(7).ColumnWidth = 14.833333333333334
$ws.Columns.Item(8).ColumnWidth = 22.666666666666668
$ws.Columns.Item(9).ColumnWidth = 25.666666666666668
$ws.Columns.Item(10).ColumnWidth = 26.5
$ws.Columns.Item(11).ColumnWidth = 26.0
$ws.Columns.Item(12).ColumnWidth = 37.666666666666664
$ws.Columns.Item(13).ColumnWidth = 2.1666666666666665

# --- Restore the active selection to where the user was editing ---
$ws.Range("D11").Select()
